# Prompts_Utterances.xlsx edit script
# - Adds a new "in progress" sheet (after "official") with a header row
# - Appends period punctuation to a batch of existing utterance strings
# - Appends 7 new rows (85-91) of utterance data to the "official" sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("official")

# ---------------------------------------------------------------------------
# 1) Punctuation touch-ups on existing utterances (sharedStrings edits)
# ---------------------------------------------------------------------------
$textEdits = @{
    "C24" = "I have no idea."
    "C25" = "I don't think so."
    "C26" = "I have no idea."
    "C27" = "I hope so."
    "C28" = "I'm sorry."
    "C29" = "I'm sorry to hear that."
    "C30" = "Thank you."
    "C31" = "I understand that."
    "C32" = "I went through something similar."
    "C33" = "I am not sure."
    "C34" = "I can't hear you."
    "C35" = "I don't understand."
    "C39" = "It would help if you could keep the questions short and simple."
    "C43" = "I don't know what that is."
    "C44" = "I don't know who that is."
    "C45" = "I am not going to talk about that."
    "C46" = "I am not here to talk about that."
    "C48" = "I don't have enough information to talk about that."
    "C50" = "That is a great question, but unfortunately I never recorded an answer to that."
    "C53" = "Unfortunately, I was never asked that question."
    "C74" = "I am stored in this machine, I can wait for your questions all day."
    "C80" = "I have answered that question before. You can ask me something else."
    "C82" = "That's a great question, but unfortunately, I don't have an answer for that right now."
    "C84" = "That is a great question but unfortunately, I don't have an answer right now."
}

foreach ($addr in $textEdits.Keys) {
    $ws.Range($addr).Value = $textEdits[$addr]
}

# ---------------------------------------------------------------------------
# 2) New rows of utterance/prompt data appended to the "official" sheet
# ---------------------------------------------------------------------------
$newRows = @(
    @("_FEEDBACK_",  "clintanderson", "Hello"),
    @("_FEEDBACK_",  "clintanderson", "It's nice to meet you."),
    @("_FEEDBACK_",  "clintanderson", "Can you please refrain from the use of profanity?"),
    @("_OFF_TOPIC_", "clintanderson", "I don't really have an opinion on that."),
    @("_FEEDBACK_",  "clintanderson", "I see you didn't ask me anything there."),
    @("_FEEDBACK_",  "clintanderson", "You've got to ask me something to get a response."),
    @("_OFF_TOPIC_", "clintanderson", "No matter what you heard about having a mouth like a sailor, it will hurt your career.")
)

$r = 85
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 3) Add the "in progress" worksheet after "official"
# ---------------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($null, $ws)
$newSheet.Name = "in progress"

$newSheet.Range("A1").Value = "Situation"
$newSheet.Range("B1").Value = "Mentor"
$newSheet.Range("C1").Value = "Utterance/Prompt"
